# Updates referee stats on "Главные" and "Линейные" worksheets and refreshes
# the "as_of_utc" timestamp (column AA) for every data row (2-26) on both sheets.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-20 03:04:15"

# ---------------------------------------------------------------------------
# Sheet "Главные"
# ---------------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Главные")

$wsMain.Range("C2").Value = 27
$wsMain.Range("D2").Value = 595
$wsMain.Range("E2").Value = 245
$wsMain.Range("F2").Value = 350
$wsMain.Range("G2").Value = 22.04
$wsMain.Range("H2").Value = 9.07
$wsMain.Range("I2").Value = 12.96
$wsMain.Range("J2").Value = 105
$wsMain.Range("K2").Value = 130

$wsMain.Range("C11").Value = 19
$wsMain.Range("D11").Value = 490
$wsMain.Range("E11").Value = 226
$wsMain.Range("F11").Value = 264
$wsMain.Range("G11").Value = 25.79
$wsMain.Range("H11").Value = 11.89
$wsMain.Range("I11").Value = 13.89
$wsMain.Range("J11").Value = 98
$wsMain.Range("K11").Value = 87
$wsMain.Range("V11").Value = 14
$wsMain.Range("W11").Value = 12

$wsMain.Range("C24").Value = 27
$wsMain.Range("D24").Value = 445
$wsMain.Range("E24").Value = 199
$wsMain.Range("F24").Value = 246
$wsMain.Range("G24").Value = 16.48
$wsMain.Range("H24").Value = 7.37
$wsMain.Range("I24").Value = 9.11
$wsMain.Range("J24").Value = 97
$wsMain.Range("K24").Value = 108

$wsMain.Range("C25").Value = 27
$wsMain.Range("D25").Value = 449
$wsMain.Range("E25").Value = 216
$wsMain.Range("F25").Value = 233
$wsMain.Range("G25").Value = 16.63
$wsMain.Range("I25").Value = 8.63
$wsMain.Range("J25").Value = 103
$wsMain.Range("K25").Value = 109
$wsMain.Range("V25").Value = 8
$wsMain.Range("W25").Value = 14

for ($r = 2; $r -le 26; $r++) {
    $wsMain.Range("AA$r").Value = $newTimestamp
}

# ---------------------------------------------------------------------------
# Sheet "Линейные"
# ---------------------------------------------------------------------------
$wsLin = $wb.Worksheets.Item("Линейные")

$wsLin.Range("C18").Value = 28
$wsLin.Range("D18").Value = 463
$wsLin.Range("E18").Value = 220
$wsLin.Range("F18").Value = 243
$wsLin.Range("G18").Value = 16.54
$wsLin.Range("H18").Value = 7.86
$wsLin.Range("I18").Value = 8.68
$wsLin.Range("J18").Value = 105
$wsLin.Range("K18").Value = 104

$wsLin.Range("C20").Value = 17
$wsLin.Range("D20").Value = 276
$wsLin.Range("E20").Value = 141
$wsLin.Range("F20").Value = 135
$wsLin.Range("G20").Value = 16.24
$wsLin.Range("H20").Value = 8.29
$wsLin.Range("I20").Value = 7.94
$wsLin.Range("J20").Value = 68
$wsLin.Range("K20").Value = 65

$wsLin.Range("C21").Value = 29
$wsLin.Range("D21").Value = 573
$wsLin.Range("E21").Value = 240
$wsLin.Range("F21").Value = 333
$wsLin.Range("G21").Value = 19.76
$wsLin.Range("H21").Value = 8.28
$wsLin.Range("I21").Value = 11.48
$wsLin.Range("J21").Value = 110
$wsLin.Range("K21").Value = 134
$wsLin.Range("V21").Value = 14
$wsLin.Range("W21").Value = 10

$wsLin.Range("C24").Value = 28
$wsLin.Range("D24").Value = 504
$wsLin.Range("E24").Value = 201
$wsLin.Range("F24").Value = 303
$wsLin.Range("G24").Value = 18
$wsLin.Range("H24").Value = 7.18
$wsLin.Range("I24").Value = 10.82
$wsLin.Range("J24").Value = 88
$wsLin.Range("K24").Value = 114
$wsLin.Range("V24").Value = 14
$wsLin.Range("W24").Value = 18

for ($r = 2; $r -le 26; $r++) {
    $wsLin.Range("AA$r").Value = $newTimestamp
}
